$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "ListFactory"
$ws.Range("Q2").Value = "K"
$ws.Range("Q3").Value = "M"

$ws.Range("Q1").Select()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
